$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShipmentTracking")

# Update row 2 data to reflect new shipment status
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "test24"
$ws.Range("D2").Value = "Reached Nearest Hub"
$ws.Range("E2").Value = "Out For Delivery"
$ws.Range("I2").Value = $false

# Move active selection to H2
$ws.Range("H2").Select()
